$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 359). All of these were bumped by one day
# (45181 -> 45182, i.e. 2023-09-12 -> 2023-09-13).
for ($r = 2; $r -le 359; $r++) {
    $ws.Cells.Item($r, 3).Value = 45182
}
